# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 should look like the rest of row 1 (bold,
# bordered, centered-top style). Copy the formatting from the last
# existing header cell (AC1) onto the new header range, then set text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row gets the same team record values.
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 83
    $ws.Cells.Item($r, 31).Value = 79
    $ws.Cells.Item($r, 32).Value = 0
}
